$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1928.3
$ws.Range("I86").Value = 1475
$ws.Range("J86").Value = 2230.5
$ws.Range("K86").Value = 1475
$ws.Range("L86").Value = 2230.5
$ws.Range("M86").Value = -352
$ws.Range("N86").Value = -4476.5
$ws.Range("H89").Value = 1928.3
$ws.Range("I89").Value = 1475
$ws.Range("J89").Value = 2230.5
$ws.Range("K89").Value = 7375
$ws.Range("L89").Value = 11152.5
$ws.Range("M89").Value = -1759
$ws.Range("N89").Value = -22384.5
$ws.Range("H116").Value = 1570489.1
$ws.Range("I116").Value = 2466340
$ws.Range("J116").Value = 2750.25
$ws.Range("K116").Value = 2466340
$ws.Range("L116").Value = 2750.25
$ws.Range("M116").Value = -2462898
$ws.Range("N116").Value = -9634.25
$ws.Range("H129").Value = 1026.3429
$ws.Range("I129").Value = 345
$ws.Range("J129").Value = 1114.258
$ws.Range("K129").Value = 1035
$ws.Range("L129").Value = 3342.774
$ws.Range("M129").Value = 3965
$ws.Range("N129").Value = -13342.774
$ws.Range("H137").Value = 1226.48
$ws.Range("I137").Value = 1164.5555
$ws.Range("J137").Value = 1385.7142
$ws.Range("K137").Value = 3493.6665
$ws.Range("L137").Value = 4157.142599999999
$ws.Range("M137").Value = -943.6664999999998
$ws.Range("N137").Value = -9257.142599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18874.527
$ws.Range("I32").Value = 3165.6826
$ws.Range("J32").Value = 128836.445
$ws.Range("K32").Value = 3165.6826
$ws.Range("L32").Value = 128836.445
$ws.Range("M32").Value = -2878.6826
$ws.Range("N32").Value = -129410.445
$ws.Range("H36").Value = 50000
$ws.Range("I36").Value = 50000
$ws.Range("K36").Value = 50000
$ws.Range("M36").Value = -49654
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H61").Value = 3836.0645
$ws.Range("I61").Value = 3078.25
$ws.Range("K61").Value = 3078.25
$ws.Range("M61").Value = -2866.25
$ws.Range("H88").Value = 1142.5
$ws.Range("I88").Value = 1185.1428
$ws.Range("J88").Value = 1099.8572
$ws.Range("K88").Value = 1185.1428
$ws.Range("L88").Value = 1099.8572
$ws.Range("M88").Value = -779.1428000000001
$ws.Range("N88").Value = -1911.8572
$ws.Range("H91").Value = 1142.5
$ws.Range("I91").Value = 1185.1428
$ws.Range("J91").Value = 1099.8572
$ws.Range("K91").Value = 1185.1428
$ws.Range("L91").Value = 1099.8572
$ws.Range("M91").Value = 218.8571999999999
$ws.Range("N91").Value = -3907.8572
$ws.Range("H132").Value = 4090.4375
$ws.Range("I132").Value = 3725.7693
$ws.Range("J132").Value = 5670.6665
$ws.Range("K132").Value = 11177.3079
$ws.Range("L132").Value = 17011.9995
$ws.Range("M132").Value = -8647.3079
$ws.Range("N132").Value = -22071.9995
$ws.Range("H136").Value = 3836.0645
$ws.Range("I136").Value = 3078.25
$ws.Range("K136").Value = 9234.75
$ws.Range("M136").Value = -6684.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10554.728
$ws.Range("I86").Value = 1874.75
$ws.Range("K86").Value = 1874.75
$ws.Range("M86").Value = -751.75
$ws.Range("H89").Value = 10554.728
$ws.Range("I89").Value = 1874.75
$ws.Range("K89").Value = 9373.75
$ws.Range("M89").Value = -3757.75
$ws.Range("H94").Value = 921.36365
$ws.Range("I94").Value = 655
$ws.Range("J94").Value = 1631.6666
$ws.Range("K94").Value = 655
$ws.Range("L94").Value = 1631.6666
$ws.Range("M94").Value = -204
$ws.Range("N94").Value = -2533.6666
$ws.Range("H134").Value = 3609.1765
$ws.Range("I134").Value = 2648.5
$ws.Range("K134").Value = 7945.5
$ws.Range("M134").Value = -5410.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1088.2354
$ws.Range("I16").Value = 908.3333
$ws.Range("J16").Value = 1520
$ws.Range("K16").Value = 908.3333
$ws.Range("L16").Value = 1520
$ws.Range("M16").Value = -621.3333
$ws.Range("N16").Value = -2094
$ws.Range("H31").Value = 6530.9355
$ws.Range("I31").Value = 2614.7
$ws.Range("J31").Value = 8395.809999999999
$ws.Range("K31").Value = 2614.7
$ws.Range("L31").Value = 8395.809999999999
$ws.Range("M31").Value = -2319.7
$ws.Range("N31").Value = -8985.809999999999
$ws.Range("H34").Value = 6530.9355
$ws.Range("I34").Value = 2614.7
$ws.Range("J34").Value = 8395.809999999999
$ws.Range("K34").Value = 2614.7
$ws.Range("L34").Value = 8395.809999999999
$ws.Range("M34").Value = -2412.7
$ws.Range("N34").Value = -8799.809999999999
$ws.Range("H113").Value = 1088.2354
$ws.Range("I113").Value = 908.3333
$ws.Range("J113").Value = 1520
$ws.Range("K113").Value = 908.3333
$ws.Range("L113").Value = 1520
$ws.Range("M113").Value = 1261.6667
$ws.Range("N113").Value = -5860
$ws.Range("H122").Value = 887.6667
$ws.Range("I122").Value = 887.6667
$ws.Range("K122").Value = 2663.0001
$ws.Range("M122").Value = -213.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 40486.4
$ws.Range("J52").Value = 40486.4
$ws.Range("L52").Value = 121459.2
$ws.Range("N52").Value = -121991.2
$ws.Range("H107").Value = 356.10526
$ws.Range("I107").Value = 356.4091
$ws.Range("J107").Value = 355.6875
$ws.Range("K107").Value = 1069.2273
$ws.Range("L107").Value = 1067.0625
$ws.Range("M107").Value = 850.7727
$ws.Range("N107").Value = -4907.0625
$ws.Range("H131").Value = 7408834.5
$ws.Range("J131").Value = 8773552
$ws.Range("L131").Value = 26320656
$ws.Range("N131").Value = -26330736
$ws.Range("H138").Value = 1205
$ws.Range("I138").Value = 886
$ws.Range("J138").Value = 2800
$ws.Range("K138").Value = 2658
$ws.Range("L138").Value = 8400
$ws.Range("M138").Value = 2482
$ws.Range("N138").Value = -18680
$ws.Range("H139").Value = 3033.7693
$ws.Range("I139").Value = 3036.5833
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 9109.749899999999
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -3969.749899999999
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3913.0588
$ws.Range("I132").Value = 4522.4
$ws.Range("J132").Value = 3042.5715
$ws.Range("K132").Value = 13567.2
$ws.Range("L132").Value = 9127.7145
$ws.Range("M132").Value = -11037.2
$ws.Range("N132").Value = -14187.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 18666.666
$ws.Range("I14").Value = 50000
$ws.Range("K14").Value = 50000
$ws.Range("M14").Value = -49832
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

Write-Output "Applied all Titan_Profits updates"